# Update "想去人数" (number of attendees wanting to go) figures for a few
# entries on both the "展览" and "全部类型" worksheets, reflecting the
# latest scrape of the source data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 3340
    $ws.Range("F5").Value = 1408
    $ws.Range("F6").Value = 18
}
